# Generate Report for Handback
# Applies the "handback" status update to the localization-status workbook:
#  - Overview / per-language sheets: status flips from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: fill in the "Latest Target File", "Latest Handback
#    File" and "Latest Handback DateTime" columns (I/J/K) for both data rows,
#    add hyperlinks on the new "Latest Target File" cells, and widen a few
#    columns that now hold longer content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared by Overview!E2:F3 and the two language sheets' Status column C)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: Latest Target File (I), Latest Handback File (J) and
#    Latest Handback DateTime (K) for rows 2 and 3.
# ---------------------------------------------------------------------------
$zhMd1 = "63745159-247d-49d2-9039-359e25082e43.md"
$zhMd2 = "b484d619-11f1-46af-8cac-4da6165ab831.md"
$zhUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a8b631c312d07469872bc42c2f678d4ab7da49c/e2e/63745159-247d-49d2-9039-359e25082e43.md"
$zhUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a8b631c312d07469872bc42c2f678d4ab7da49c/e2e/b484d619-11f1-46af-8cac-4da6165ab831.md"

$zh.Range("J2").Value = "63745159-247d-49d2-9039-359e25082e43.50ef2da33002855dd407f94bb5af1e67736546c2.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-06 03:40:15"

$zh.Range("J3").Value = "b484d619-11f1-46af-8cac-4da6165ab831.8868e297c8e153731a71001790029906cd8cb4b7.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-06 03:40:15"

# Rebuild the hyperlink list so I2/I3 pick up hyperlinks to the source .md
# files (same targets as A2/A3) alongside the pre-existing A2/A3 links.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhUrl1, "", "", $zhMd1)
$zh.Range("I2").Value = $zhMd1
$zh.Hyperlinks.Add($zh.Range("I2"), $zhUrl1, "", "", $zhMd1)
$zh.Hyperlinks.Add($zh.Range("A3"), $zhUrl2, "", "", $zhMd2)
$zh.Range("I3").Value = $zhMd2
$zh.Hyperlinks.Add($zh.Range("I3"), $zhUrl2, "", "", $zhMd2)

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment, with its own handback datetime.
# ---------------------------------------------------------------------------
$deMd1 = "63745159-247d-49d2-9039-359e25082e43.md"
$deMd2 = "b484d619-11f1-46af-8cac-4da6165ab831.md"
$deUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a8b631c312d07469872bc42c2f678d4ab7da49c/e2e/63745159-247d-49d2-9039-359e25082e43.md"
$deUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a8b631c312d07469872bc42c2f678d4ab7da49c/e2e/b484d619-11f1-46af-8cac-4da6165ab831.md"

$de.Range("J2").Value = "63745159-247d-49d2-9039-359e25082e43.50ef2da33002855dd407f94bb5af1e67736546c2.de-de.xlf"
$de.Range("K2").Value = "2016-09-06 03:40:31"

$de.Range("J3").Value = "b484d619-11f1-46af-8cac-4da6165ab831.8868e297c8e153731a71001790029906cd8cb4b7.de-de.xlf"
$de.Range("K3").Value = "2016-09-06 03:40:31"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deUrl1, "", "", $deMd1)
$de.Range("I2").Value = $deMd1
$de.Hyperlinks.Add($de.Range("I2"), $deUrl1, "", "", $deMd1)
$de.Hyperlinks.Add($de.Range("A3"), $deUrl2, "", "", $deMd2)
$de.Range("I3").Value = $deMd2
$de.Hyperlinks.Add($de.Range("I3"), $deUrl2, "", "", $deMd2)

# ---------------------------------------------------------------------------
# 4. Column width adjustments (columns now show longer status/file text).
# ---------------------------------------------------------------------------
$ov.Range("E1").ColumnWidth = 29.9777050018311
$ov.Range("F1").ColumnWidth = 29.9777050018311

$zh.Range("C1").ColumnWidth = 29.9777050018311
$zh.Range("I1").ColumnWidth = 40
$zh.Range("J1").ColumnWidth = 40

$de.Range("C1").ColumnWidth = 29.9777050018311
$de.Range("I1").ColumnWidth = 40
$de.Range("J1").ColumnWidth = 40
